$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the existing data rows (use Value2 for reliable reads of cell
# contents, including blanks) before shifting them down.
$row2A = $ws.Range("A2").Value2
$row3A = $ws.Range("A3").Value2
$row3B = $ws.Range("B3").Value2
$row3C = $ws.Range("C3").Value2

# Shift old row 3 (Unassigned | Unassigned | Unassigned) down to row 4.
$ws.Range("A4").Value = $row3A
$ws.Range("B4").Value = $row3B
$ws.Range("C4").Value = $row3C

# Clear the old row 3 cells so the move doesn't leave duplicate data behind.
$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Shift old row 2 (Cololabis saira) down to row 3.
$ws.Range("A3").Value = $row2A

# Write the new species name into row 2.
$ws.Range("A2").Value = "Fundulus heteroclitus or majalis"

# Append two new rows with additional species names.
$ws.Range("A5").Value = "Mareca americana"
$ws.Range("A6").Value = "Myrophis vafer"
